$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new contingency "line" entries (line7, line8) are inserted into the
# shared-string table right after line6 / before the extrN block. Since the
# underlying string index for B8/B9 is unchanged by the edit, their displayed
# label now resolves to line7/line8 instead of extr1/extr2, and every row
# that used to read extrN (N>=3) now reads extr(N-2). Two brand-new rows
# (16,17) are appended holding extr7/extr8.

# Row 8 (was extr1 -> line7): C,D values change, E stays False
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $false

# Row 9 (was extr2 -> line8): C changes, D stays, E flips to True
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 5).Value = $true

# Row 10 (was extr3 -> extr1): C,D change, E flips to True
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# Row 11 (was extr4 -> extr2): C,D change, E flips to True
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# Row 12 (was extr5 -> extr3): C changes, D/E stay
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(12, 3).Value = 10

# Row 13 (was extr6 -> extr4): D changes, E flips to False
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false

# Row 14 (was extr7 -> extr5): C,D change, E flips to False
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $false

# Row 15 (was extr8 -> extr6): C,D change, E flips to True
$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

# New row 16: extr7
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $true

# New row 17: extr8
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false

# Column A on the new rows carries the same bold/centered/bordered look as
# the rest of column A; copy that formatting down from row 15.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
